# Apply the 2018 "Doing Updates for Financials" data refresh to the NGG yearly financials sheet.
# Updates the numeric values (columns D:J) for each changed financial statement line item.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Total Revenue
$ws.Range("D8").Value = 19890700
$ws.Range("E8").Value = 19610300
$ws.Range("F8").Value = 17232500
$ws.Range("G8").Value = 17421700
$ws.Range("H8").Value = 19315500
$ws.Range("I8").Value = 18728600
$ws.Range("J8").Value = 18041200

# Row 9: Cost of Revenue
$ws.Range("D9").Value = 7495900
$ws.Range("E9").Value = 6888100
$ws.Range("F9").Value = 6251600
$ws.Range("G9").Value = 7082400
$ws.Range("H9").Value = 7759300
$ws.Range("I9").Value = 6560700
$ws.Range("J9").Value = 7098100

# Row 10: Gross Profit
$ws.Range("D10").Value = 12394900
$ws.Range("E10").Value = 12722200
$ws.Range("F10").Value = 10981000
$ws.Range("G10").Value = 10339300
$ws.Range("H10").Value = 11556200
$ws.Range("I10").Value = 12167900
$ws.Range("J10").Value = 10943200

# Row 12: Research Development
$ws.Range("D12").Value = 17000
$ws.Range("E12").Value = 18300
$ws.Range("F12").Value = 24800
$ws.Range("G12").Value = 20900
$ws.Range("H12").Value = 15700
$ws.Range("I12").Value = 19600
$ws.Range("J12").Value = 19600

# Row 14: Non Recurring
$ws.Range("D14").Value = -33900
$ws.Range("E14").Value = 825600
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 170900
$ws.Range("I14").Value = -232200
$ws.Range("J14").Value = 397800

# Row 17: Total Operating Expenses
$ws.Range("D17").Value = 15334800
$ws.Range("E17").Value = 15426100
$ws.Range("F17").Value = 13026100
$ws.Range("G17").Value = 13743500
$ws.Range("H17").Value = 14443900
$ws.Range("I17").Value = 13838700
$ws.Range("J17").Value = 13430500

# Row 18: Operating Income or Loss
$ws.Range("D18").Value = 4556000
$ws.Range("E18").Value = 4184200
$ws.Range("F18").Value = 4206400
$ws.Range("G18").Value = 3678200
$ws.Range("H18").Value = 4871600
$ws.Range("I18").Value = 4889900
$ws.Range("J18").Value = 4610700

# Row 20: Total Other Income/Expenses Net
$ws.Range("D20").Value = 266100
$ws.Range("E20").Value = -181300
$ws.Range("F20").Value = -293500
$ws.Range("G20").Value = -95200
$ws.Range("H20").Value = 45700
$ws.Range("I20").Value = -56100
$ws.Range("J20").Value = -114800

# Row 21: Earnings Before Interest And Taxes
$ws.Range("D21").Value = 6840000
$ws.Range("E21").Value = 5956200
$ws.Range("F21").Value = 5642000
$ws.Range("G21").Value = 5168300
$ws.Range("H21").Value = 6786100
$ws.Range("I21").Value = 6628800
$ws.Range("J21").Value = 6186800

# Row 22: Interest Expense
$ws.Range("D22").Value = 1290000
$ws.Range("E22").Value = 1154300
$ws.Range("F22").Value = 875200
$ws.Range("G22").Value = 1026500
$ws.Range("H22").Value = 1333000
$ws.Range("I22").Value = 1297800
$ws.Range("J22").Value = 1389100

# Row 23: Income Before Tax
$ws.Range("D23").Value = 3532100
$ws.Range("E23").Value = 2848600
$ws.Range("F23").Value = 3037700
$ws.Range("G23").Value = 2556400
$ws.Range("H23").Value = 3584200
$ws.Range("I23").Value = 3536000
$ws.Range("J23").Value = 3106900

# Row 24: Income Tax Expense
$ws.Range("D24").Value = -1153000
$ws.Range("E24").Value = 487800
$ws.Range("F24").Value = 556900
$ws.Range("G24").Value = 609100
$ws.Range("H24").Value = 370400
$ws.Range("I24").Value = 726500
$ws.Range("J24").Value = 603900

# Row 26: Income After Tax
$ws.Range("D26").Value = 4685100
$ws.Range("E26").Value = 2360800
$ws.Range("F26").Value = 2480800
$ws.Range("G26").Value = 1947300
$ws.Range("H26").Value = 3213800
$ws.Range("I26").Value = 2809500
$ws.Range("J26").Value = 2503000

# Row 27: Net Income From Continuing Ops
$ws.Range("D27").Value = 4683800
$ws.Range("E27").Value = 2360800
$ws.Range("F27").Value = 2479500
$ws.Range("G27").Value = 1960400
$ws.Range("H27").Value = 3229500
$ws.Range("I27").Value = 2808200
$ws.Range("J27").Value = 2500400

# Row 29: Discontinued Operations
$ws.Range("D29").Value = -53500
$ws.Range("E29").Value = 7806300
$ws.Range("F29").Value = 900000
$ws.Range("G29").Value = 673000

# Row 32: Other Items
$ws.Range("D32").Value = -266100
$ws.Range("E32").Value = 181300
$ws.Range("F32").Value = 293500
$ws.Range("G32").Value = 95200
$ws.Range("H32").Value = -45700
$ws.Range("I32").Value = 56100
$ws.Range("J32").Value = 114800

# Row 33: Net Income
$ws.Range("D33").Value = 4630300
$ws.Range("E33").Value = 10167100
$ws.Range("F33").Value = 3379500
$ws.Range("G33").Value = 2633400
$ws.Range("H33").Value = 3229500
$ws.Range("I33").Value = 2808200
$ws.Range("J33").Value = 2500400

# Row 35: Net Income Applicable To Common Shares
$ws.Range("D35").Value = 4630300
$ws.Range("E35").Value = 10167100
$ws.Range("F35").Value = 3379500
$ws.Range("G35").Value = 2633400
$ws.Range("H35").Value = 3229500
$ws.Range("I35").Value = 2808200
$ws.Range("J35").Value = 2500400

# Row 41: Cash And Cash Equivalents
$ws.Range("D41").Value = 429100
$ws.Range("E41").Value = 1485600
$ws.Range("F41").Value = 165600
$ws.Range("G41").Value = 155200
$ws.Range("H41").Value = 461700
$ws.Range("I41").Value = 875200
$ws.Range("J41").Value = 433000

# Row 42: Short Term Investments
$ws.Range("D42").Value = 3005100
$ws.Range("E42").Value = 9693600
$ws.Range("F42").Value = 2544700
$ws.Range("G42").Value = 1606900
$ws.Range("H42").Value = 3542500
$ws.Range("I42").Value = 12876100
$ws.Range("J42").Value = 3118600

# Row 43: Net Receivables
$ws.Range("D43").Value = 4008100
$ws.Range("E43").Value = 5381600
$ws.Range("F43").Value = 4284700
$ws.Range("G43").Value = 3974200
$ws.Range("H43").Value = 3399000
$ws.Range("I43").Value = 4906800
$ws.Range("J43").Value = 1269100

# Row 44: Inventory
$ws.Range("D44").Value = 444800
$ws.Range("E44").Value = 525600
$ws.Range("F44").Value = 570000
$ws.Range("G44").Value = 443500
$ws.Range("H44").Value = 349600
$ws.Range("I44").Value = 390000
$ws.Range("J44").Value = 490400

# Row 45: Other Current Assets
$ws.Range("D45").Value = 826900
$ws.Range("E45").Value = 688700
$ws.Range("F45").Value = 667800
$ws.Range("G45").Value = 1686500
$ws.Range("H45").Value = 2015200
$ws.Range("I45").Value = 2264300
$ws.Range("J45").Value = 1715200

# Row 46: Total Current Assets
$ws.Range("D46").Value = 8714100
$ws.Range("E46").Value = 17704700
$ws.Range("F46").Value = 8232800
$ws.Range("G46").Value = 7866300
$ws.Range("H46").Value = 9768000
$ws.Range("I46").Value = 12490100
$ws.Range("J46").Value = 7026300

# Row 47: Long Term Investments
$ws.Range("D47").Value = 4047300
$ws.Range("E47").Value = 4210300
$ws.Range("F47").Value = 1194700
$ws.Range("G47").Value = 896100
$ws.Range("H47").Value = 871300
$ws.Range("I47").Value = 979500
$ws.Range("J47").Value = 842600

# Row 48: Property Plant and Equipment
$ws.Range("D48").Value = 51980700
$ws.Range("E48").Value = 51944100
$ws.Range("F48").Value = 56560100
$ws.Range("G48").Value = 53115400
$ws.Range("H48").Value = 48492900
$ws.Range("I48").Value = 27038300
$ws.Range("J48").Value = 43956600

# Row 49: Goodwill
$ws.Range("D49").Value = 8273200
$ws.Range("E49").Value = 9155000
$ws.Range("F49").Value = 8089300
$ws.Range("G49").Value = 7756700
$ws.Range("H49").Value = 6864600
$ws.Range("I49").Value = 8094500
$ws.Range("J49").Value = 6941500

# Row 52: Other Assets
$ws.Range("D52").Value = 3661200
$ws.Range("E52").Value = 2929500
$ws.Range("F52").Value = 2791200
$ws.Range("G52").Value = 2218600
$ws.Range("H52").Value = 2328200
$ws.Range("I52").Value = 2964700
$ws.Range("J52").Value = 2972500

# Row 54: Total Assets
$ws.Range("D54").Value = 76676500
$ws.Range("E54").Value = 85875800
$ws.Range("F54").Value = 76868200
$ws.Range("G54").Value = 71853100
$ws.Range("H54").Value = 68325000
$ws.Range("I54").Value = 71352300
$ws.Range("J54").Value = 61739500

# Row 57: Accounts Payable
$ws.Range("D57").Value = 2578600
$ws.Range("E57").Value = 2784700
$ws.Range("F57").Value = 2658200
$ws.Range("G57").Value = 2673800
$ws.Range("H57").Value = 2533000
$ws.Range("I57").Value = 5303300
$ws.Range("J57").Value = 1995600

# Row 58: Short/Current Long Term Debt
$ws.Range("D58").Value = 5800300
$ws.Range("E58").Value = 7168500
$ws.Range("F58").Value = 4709900
$ws.Range("G58").Value = 3949500
$ws.Range("H58").Value = 4579400
$ws.Range("I58").Value = 5084200
$ws.Range("J58").Value = 3250300

# Row 59: Other Current Liabilities
$ws.Range("D59").Value = 2964700
$ws.Range("E59").Value = 3877700
$ws.Range("F59").Value = 2702500
$ws.Range("G59").Value = 2994700
$ws.Range("H59").Value = 2449500
$ws.Range("I59").Value = 2561700
$ws.Range("J59").Value = 2585100

# Row 60: Total Current Liabilities
$ws.Range("D60").Value = 11343600
$ws.Range("E60").Value = 13709600
$ws.Range("F60").Value = 10070600
$ws.Range("G60").Value = 9618000
$ws.Range("H60").Value = 9561900
$ws.Range("I60").Value = 9710600
$ws.Range("J60").Value = 7831100

# Row 61: Long Term Debt
$ws.Range("D61").Value = 28927000
$ws.Range("E61").Value = 30184300
$ws.Range("F61").Value = 32259500
$ws.Range("G61").Value = 29845200
$ws.Range("H61").Value = 29267400
$ws.Range("I61").Value = 32147300
$ws.Range("J61").Value = 26781400

# Row 62: Other Liabilities
$ws.Range("D62").Value = 11822300
$ws.Range("E62").Value = 15394800
$ws.Range("F62").Value = 16845200
$ws.Range("G62").Value = 16772100
$ws.Range("H62").Value = 13949600
$ws.Range("I62").Value = 16146100
$ws.Range("J62").Value = 15067400

# Row 66: Total Liabilities
$ws.Range("D66").Value = 52113700
$ws.Range("E66").Value = 59309600
$ws.Range("F66").Value = 59188300
$ws.Range("G66").Value = 56251000
$ws.Range("H66").Value = 52789300
$ws.Range("I66").Value = 58010500
$ws.Range("J66").Value = 49689000

# Row 72: Retained Earnings
$ws.Range("D72").Value = 21459800
$ws.Range("E72").Value = 22741900
$ws.Range("F72").Value = 14554800
$ws.Range("G72").Value = 12683100
$ws.Range("H72").Value = 12715700
$ws.Range("I72").Value = 10324900
$ws.Range("J72").Value = 9196700

# Row 76: Total Stockholder Equity
$ws.Range("D76").Value = 24562800
$ws.Range("E76").Value = 26566200
$ws.Range("F76").Value = 17679900
$ws.Range("G76").Value = 15602200
$ws.Range("H76").Value = 15535600
$ws.Range("I76").Value = 13341800
$ws.Range("J76").Value = 12050500

# Row 81: Net Income
$ws.Range("D81").Value = 4630300
$ws.Range("E81").Value = 10167100
$ws.Range("F81").Value = 3379500
$ws.Range("G81").Value = 2633400
$ws.Range("H81").Value = 3229500
$ws.Range("I81").Value = 2808200
$ws.Range("J81").Value = 2500400

# Row 83: Depreciation
$ws.Range("D83").Value = 1995600
$ws.Range("E83").Value = 1931700
$ws.Range("F83").Value = 1710000
$ws.Range("G83").Value = 1567800
$ws.Range("H83").Value = 1848200
$ws.Range("I83").Value = 1775200
$ws.Range("J83").Value = 1672100

# Row 89: Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 5873300
$ws.Range("E89").Value = 6839800
$ws.Range("F89").Value = 7001500
$ws.Range("G89").Value = 6530700
$ws.Range("H89").Value = 5242000
$ws.Range("I89").Value = 4891200
$ws.Range("J89").Value = 5514600

# Row 91: Capital Expenditures
$ws.Range("D91").Value = -4875500
$ws.Range("E91").Value = -4299000
$ws.Range("F91").Value = -3723800
$ws.Range("G91").Value = -3362500
$ws.Range("H91").Value = -3839900
$ws.Range("I91").Value = -4192100
$ws.Range("J91").Value = -4104700

# Row 94: Total Cash Flows From Investing Activities
$ws.Range("D94").Value = 2917700
$ws.Range("E94").Value = -5626800
$ws.Range("F94").Value = -5264200
$ws.Range("G94").Value = -2609900
$ws.Range("H94").Value = -1734700
$ws.Range("I94").Value = -7995400
$ws.Range("J94").Value = -3092500

# Row 96: Dividends Paid
$ws.Range("D96").Value = -5852400
$ws.Range("E96").Value = -1908200
$ws.Range("F96").Value = -1743900
$ws.Range("G96").Value = -1657800
$ws.Range("H96").Value = -1381300
$ws.Range("I96").Value = -1056500
$ws.Range("J96").Value = -1312100

# Row 100: Total Cash Flows From Financing Activities
$ws.Range("D100").Value = -9843600
$ws.Range("E100").Value = 90000
$ws.Range("F100").Value = -1732100
$ws.Range("G100").Value = -4242900
$ws.Range("H100").Value = -3876400
$ws.Range("I100").Value = 3541200
$ws.Range("J100").Value = -2478200

# Row 101: Effect Of Exchange Rate Changes 
$ws.Range("D101").Value = -3900
$ws.Range("E101").Value = 20900
$ws.Range("F101").Value = 5200
$ws.Range("G101").Value = 31300
$ws.Range("H101").Value = -33900
$ws.Range("I101").Value = 18300
$ws.Range("J101").Value = 0

# Row 102: Change In Cash and Cash Equivalents 
$ws.Range("D102").Value = -1056500
$ws.Range("E102").Value = 1323900
$ws.Range("F102").Value = 10400
$ws.Range("G102").Value = -290900
$ws.Range("H102").Value = -403000
$ws.Range("I102").Value = 455200
$ws.Range("J102").Value = -56100
